$p = $ppt.ActivePresentation

# Slide 2: "Small but complete programming language..." paragraph
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$shp2.TextFrame.TextRange.Paragraphs(1).Text = "Small but complete programming language with constructs similar to those found in Ada and C-based languages such as Java, Kotlin, and C#."

# Slide 4: "identifier = letter ( letter | digit )* ." -> "identifier = letter { letter | digit } ."
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$shp4.TextFrame.TextRange.Paragraphs(3).Text = "identifier = letter { letter | digit } ."

# Slide 8: "Examples include variable declarations, ..." -> add "constant declarations, "
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$shp8.TextFrame.TextRange.Paragraphs(3).Text = "Examples include constant declarations, variable declarations, type declarations, and subprogram declarations."

# Slide 19: merge the two runs of the "name := ""Caleb"";      // length = 5" paragraph into one
$s19 = $p.Slides.Item(19)
$shp19 = $s19.Shapes.Item(2)
$shp19.TextFrame.TextRange.Paragraphs(10).Text = "name := ""Caleb"";      // length = 5"
